# airdas_check bug fix, docs
#
# 1. "Checks performe" + _GoBack bookmark + "d by the function " -> merge into
#    a single run "Checks performed by the function " (bookmark dropped here).
# 2. Latitude bullet: "(inclusive)" -> "(inclusive; NA values are ignored)",
#    written as three separate runs.
# 3. Longitude bullet: "(inclusive)" -> "(inclusive; NA values are ignored)",
#    written as three separate runs, and the _GoBack bookmark is re-created
#    at the end of this paragraph (that's where Word leaves it after the
#    last edit of the session).

$d = $word.ActiveDocument

function Set-RunsXml($rng, $pPr, $innerRuns) {
    $xml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p w:rsidR="00FF48C7" w:rsidRDefault="00FF48C7" w:rsidP="00CD5B51">$pPr$innerRuns</w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
    $rng.Text = ""
    $rng.InsertXML($xml)
}

$listPPr = "<w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""1""/></w:numPr><w:spacing w:after=""0"" w:line=""240"" w:lineRule=""auto""/></w:pPr>"

# --- 1. merge the split "Checks performe" / "d by the function " runs, drop the bookmark ---
$rng = $d.Content
$rng.Find.Execute("Checks performed by the function ") | Out-Null
$rng.Text = ""
$rng.InsertXML(@"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p><w:r><w:t xml:space="preserve">Checks performed by the function </w:t></w:r></w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@)

# --- 2. Latitude bullet: split "(inclusive)" into 3 runs ---
$rng = $d.Content
$rng.Find.Execute("Latitude values are between -90 and 90 (inclusive)") | Out-Null
$latRuns = "<w:r><w:t>Latitude values are between -90 and 90 (inclusive</w:t></w:r><w:r><w:t>; NA values are ignored</w:t></w:r><w:r><w:t>)</w:t></w:r>"
Set-RunsXml $rng $listPPr $latRuns

# --- 3. Longitude bullet: split "(inclusive)" into 3 runs, re-add _GoBack bookmark at the end ---
$rng = $d.Content
$rng.Find.Execute("Longitude values are between -180 and 180 (inclusive)") | Out-Null
$lonRuns = "<w:r><w:t>Longitude values are between -180 and 180 (inclusive</w:t></w:r><w:r><w:t>; NA values are ignored</w:t></w:r><w:r><w:t>)</w:t></w:r><w:bookmarkStart w:id=""0"" w:name=""_GoBack""/><w:bookmarkEnd w:id=""0""/>"
Set-RunsXml $rng $listPPr $lonRuns
